# Generate Report for Handback
#
# Row 7 in both the "zh-cn" and "de-de" worksheets corresponds to the
# 6c927b7a-0007-4e1f-a8d5-0408416c54cf file, which has now received a
# handback. Populate the "Latest Target File" (I), "Latest Handback File"
# (J), "Latest Handback DateTime" (K) and "Error Detail" (P) columns for
# that row on both language sheets, and wire up the hyperlink for the new
# "Latest Target File" entry, matching the existing hyperlink style used
# elsewhere in the sheet (underline, CornflowerBlue font colour).

$wb = $excel.ActiveWorkbook

$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f1a2b182cbef91694c5edac85ce953f93d8ce8df/e2e/6c927b7a-0007-4e1f-a8d5-0408416c54cf.md"
$targetDisplay = "6c927b7a-0007-4e1f-a8d5-0408416c54cf.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1412e407e4b0705a835fbf000f07b3c50e8e1871/e2e/6c927b7a-0007-4e1f-a8d5-0408416c54cf.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f1a2b182cbef91694c5edac85ce953f93d8ce8df/e2e/6c927b7a-0007-4e1f-a8d5-0408416c54cf.md."
$hyperlinkColor = 15570276

# --- zh-cn sheet, row 7 ---
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $targetUrl, "", "", $targetDisplay)
$wsZh.Range("I7").Font.Underline = 2
$wsZh.Range("I7").Font.Color = $hyperlinkColor

$wsZh.Range("J7").Value = "6c927b7a-0007-4e1f-a8d5-0408416c54cf.501397463eeabf571276ccd5b1ee449fd4d4aef0.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-29 21:09:57"
$wsZh.Range("P7").Value = $errorDetail

# --- de-de sheet, row 7 ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $targetUrl, "", "", $targetDisplay)
$wsDe.Range("I7").Font.Underline = 2
$wsDe.Range("I7").Font.Color = $hyperlinkColor

$wsDe.Range("J7").Value = "6c927b7a-0007-4e1f-a8d5-0408416c54cf.501397463eeabf571276ccd5b1ee449fd4d4aef0.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-29 21:10:16"
$wsDe.Range("P7").Value = $errorDetail
